# First set of model updates for PL
# Rename the "HU" sheet to "PL" and refresh its RMSE coefficients.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HU")
$ws.Name = "PL"

$ws.Range("B2").Value = 0.93553296718914958
$ws.Range("B3").Value = 0.89086026422674336
$ws.Range("B4").Value = 0.40974616881083442
$ws.Range("B5").Value = 0.45713021107178931
$ws.Range("B6").Value = 1.0259543542924621
$ws.Range("B7").Value = 1.0118989434177876
